$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set "Execute" column (C) to "Y" for rows 2, 3, 5, 6, 7 (row 4 already "Y")
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"
$ws.Range("C5").Value = "Y"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

# Rename "Existing Liability w/Notice Number" -> "Existing Liability with Notice/Invoice Number"
$ws.Range("D2").Value = "Existing Liability with Notice/Invoice Number"
$ws.Range("D6").Value = "Existing Liability with Notice/Invoice Number"

# Update selection to C14 (matches new sheetView selection in the saved file)
$ws.Range("C14").Select()
